# Update "想去人数" (people interested) counts on the "展览" (exhibitions)
# and "全部类型" (all types) worksheets, matching the refreshed data pull.

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 25
$wsExpo.Range("F5").Value = 3282
$wsExpo.Range("F6").Value = 2111
$wsExpo.Range("F9").Value = 27
$wsExpo.Range("F10").Value = 1203
$wsExpo.Range("F12").Value = 1214
$wsExpo.Range("F13").Value = 99

# --- Sheet: 全部类型 ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 25
$wsAll.Range("F5").Value = 3282
$wsAll.Range("F6").Value = 2111
$wsAll.Range("F10").Value = 27
$wsAll.Range("F11").Value = 1203
$wsAll.Range("F13").Value = 1214
$wsAll.Range("F14").Value = 99
